$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 44584
$ws.Range("C2").Formula = "=B2/60"
